$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Sheet structure: insert "Workflows" and "Sheet2" before "FrontEnd",
#    remove "QGIS 3", hide "GG Map API".
# ---------------------------------------------------------------------------
$wsWorkflows = $wb.Worksheets.Add($wb.Worksheets.Item("FrontEnd"))
$wsWorkflows.Name = "Workflows"

$wsSheet2 = $wb.Worksheets.Add($wb.Worksheets.Item("FrontEnd"))
$wsSheet2.Name = "Sheet2"

$wb.Worksheets.Item("QGIS 3").Delete() | Out-Null
$wb.Worksheets.Item("GG Map API").Visible = $false

# ---------------------------------------------------------------------------
# 2. "Workflows" sheet content
# ---------------------------------------------------------------------------
$wsWorkflows.Range("B2").Value = "Đọc file CSV để lấy lịch sử phim đã được chiếu"
$wsWorkflows.Range("B3").Value = "Hiển thị bản đồ San Francisco tổng quan"
$wsWorkflows.Range("B4").Value = "Hiện thị toàn bộ marker các địa điểm chiếu phim"
$wsWorkflows.Range("B6").Value = 'Tiền xử lí dữ liệu "địa chỉ" (Locations) sao cho thành chuỗi truy vấn dùng cho url'
$wsWorkflows.Range("C7").Value = "_ and _ Street"
$wsWorkflows.Range("D8").Value = "20th and Folsom Streets"
$wsWorkflows.Range("D9").Value = "street=Folsom+%26+20th"

$wsWorkflows.Range("C11").Value = "https://nominatim.openstreetmap.org/search.php?street=Folsom+%26+20th&city=san+francisco&format=jsonv2"
$wsWorkflows.Hyperlinks.Add($wsWorkflows.Range("C11"), "https://nominatim.openstreetmap.org/search.php?street=Folsom+%26+20th&city=san+francisco&format=jsonv2") | Out-Null

$json1 = "[`n    {`n        ""place_id"": 68854055,`n        ""licence"": ""Data © OpenStreetMap contributors, ODbL 1.0. https://osm.org/copyright"",`n        ""osm_type"": ""node"",`n        ""osm_id"": 6384663492,`n        ""boundingbox"": [`n            ""37.758673"",`n            ""37.758773"",`n            ""-122.4148487"",`n            ""-122.4147487""`n        ],`n        ""lat"": ""37.758723"",`n        ""lon"": ""-122.4147987"",`n        ""display_name"": ""Folsom Street & 20th Street, Folsom Street, Mission District, San Francisco, California, 90103, United States"",`n        ""place_rank"": 30,`n        ""category"": ""highway"",`n        ""type"": ""bus_stop"",`n        ""importance"": 0.4001,`n        ""icon"": ""https://nominatim.openstreetmap.org/ui/mapicons/transport_bus_stop2.p.20.png""`n    }`n]"
$wsWorkflows.Range("C12").Value = $json1

$wsWorkflows.Range("C13").Value = "https://nominatim.openstreetmap.org/search.php?q=Golden+Gate+Park%2C+San+Francisco&format=jsonv2"
$wsWorkflows.Hyperlinks.Add($wsWorkflows.Range("C13"), "https://nominatim.openstreetmap.org/search.php?q=Golden+Gate+Park%2C+San+Francisco&format=jsonv2") | Out-Null

$json2 = "[`n    {`n        ""place_id"": 142363953,`n        ""licence"": ""Data © OpenStreetMap contributors, ODbL 1.0. https://osm.org/copyright"",`n        ""osm_type"": ""way"",`n        ""osm_id"": 158602261,`n        ""boundingbox"": [`n            ""37.764142"",`n            ""37.7746518"",`n            ""-122.5108673"",`n            ""-122.45318""`n        ],`n        ""lat"": ""37.769368099999994"",`n        ""lon"": ""-122.48218371117709"",`n        ""display_name"": ""Golden Gate Park, San Francisco, California, United States"",`n        ""place_rank"": 24,`n        ""category"": ""leisure"",`n        ""type"": ""park"",`n        ""importance"": 0.955740352341206`n    }`n]"
$wsWorkflows.Range("C14").Value = $json2

# Style index 2 (applyAlignment) — toggle WrapText to materialise the extra cellXfs entry.
$wsWorkflows.Range("C12").WrapText = $true
$wsWorkflows.Range("C12").WrapText = $false
$wsWorkflows.Range("C14").WrapText = $true
$wsWorkflows.Range("C14").WrapText = $false

# ---------------------------------------------------------------------------
# 3. "Sheet2" content
# ---------------------------------------------------------------------------
$webdevUrl = "https://thewebdev.info/2022/04/03/how-to-pass-variables-from-python-flask-to-javascript/#:~:text=To%20pass%20variables%20from%20Python%20Flask%20to%20JavaScript%2C%20we%20can,to%20pass%20to%20the%20template.&text=to%20call%20render_template%20with%20the,to%20pass%20to%20the%20template."
$wsSheet2.Range("B2").Value = $webdevUrl
$wsSheet2.Hyperlinks.Add($wsSheet2.Range("B2"), $webdevUrl) | Out-Null

$wsSheet2.Range("B4").Value = "<html>"
$wsSheet2.Range("B5").Value = "    <head>"
$wsSheet2.Range("B6").Value = "         <script type=""text/javascript"" {{ url_for('static', filename='app.js')}}></script>"
$wsSheet2.Range("B7").Value = "         <script type=""text/javascript"">"
$wsSheet2.Range("B8").Value = "            myVar = myFunc({{ data | tojson }})"
$wsSheet2.Range("B9").Value = "         </script>"
$wsSheet2.Range("B10").Value = "    </head>"
$wsSheet2.Range("B11").Value = "</html>"

# ---------------------------------------------------------------------------
# 4. "OpenStreetMap" sheet: add new nominatim row
# ---------------------------------------------------------------------------
$wsOsm = $wb.Worksheets.Item("OpenStreetMap")
$wsOsm.Range("B7").Value = "https://nominatim.org/release-docs/develop/api/Search/"
$wsOsm.Hyperlinks.Add($wsOsm.Range("B7"), "https://nominatim.org/release-docs/develop/api/Search/") | Out-Null

# ---------------------------------------------------------------------------
# 5. Selection / active sheet state
# ---------------------------------------------------------------------------
$wsOsm.Range("B7").Select() | Out-Null
$wsSheet2.Range("F9").Select() | Out-Null
$wsWorkflows.Range("D9").Select() | Out-Null
$wb.Worksheets.Item("Workflows").Activate() | Out-Null
